$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Google")
$ws2 = $wb.Worksheets.Item("DuckDuckGo")

# --- "Google" sheet: replace the old test keyword with the new DuckDuckGo-era
# keyword set, adding two more rows ---
$ws1.Range("A2").Value = "automation testing"
$ws1.Range("A3").Value = "quality assurance"
$ws1.Range("A4").Value = "ruby is the best programming language"

# --- "DuckDuckGo" sheet content stays the same, only the view changes ---

# --- Selections / active sheet ---
# Select A3 on Google first (leaves it as the non-active tab selection).
$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null

# Make DuckDuckGo the final active / selected tab, with B7 selected.
$ws2.Activate() | Out-Null
$ws2.Range("B7").Select() | Out-Null
